$d = $word.ActiveDocument

# 1. Merge the "List and describe 3 different deployment strategies..." run
#    with the following single-space run into one run (the trailing
#    whitespace run is removed and its content appended to the first run).
$old1 = "List and describe 3 different deployment strategies that Parts Unlimited could have used instead of their selected approach. What are the benefits? What are the drawbacks?"
$found1 = $d.Content.Find.Execute($old1, $false, $false, $false, $false, $false, $true, 1, $false, $old1, 2)

# 2. Remove the sentence about convincing CEO Steve to back the "growth
#    thesis" from the closing paragraph.
$old2 = "Their combined efforts convinced CEO Steve to back the “growth thesis,” allocate `$5 million for an Innovation Team, and endorse ongoing exploratory efforts—despite Sarah’s objections. "
$found2 = $d.Content.Find.Execute($old2, $false, $false, $false, $false, $false, $true, 1, $false, "", 2)

Write-Output "Edit1 found: $found1; Edit2 found: $found2"
